{"js": "const replacements = [\n    [\"2023-09-03 Sunday\", \"2023-09-04 Monday\"],\n    [\"33\u00d747=1551\", \"83\u00d711=913\"],\n    [\"73\u00d768=4964\", \"69\u00d714=966\"],\n    [\"62\u00d742=2604\", \"40\u00d751=2040\"],\n    [\"79\u00d795=7505\", \"25\u00d798=2450\"],\n    [\"62\u00d747=2914\", \"28\u00d780=2240\"],\n    [\"72\u00d755=3960\", \"41\u00d740=1640\"],\n    [\"44\u00d728=1232\", \"17\u00d713=221\"],\n    [\"70\u00d723=1610\", \"65\u00d758=3770\"],\n    [\"43\u00d716=688\", \"88\u00d726=2288\"],\n    [\"42\u00d797=4074\", \"92\u00d719=1748\"],\n    [\"32\u00d796=3072\", \"67\u00d725=1675\"],\n    [\"19\u00d726=494\", \"38\u00d731=1178\"],\n    [\"27\u00d783=2241\", \"98\u00d784=8232\"],\n    [\"19\u00d756=1064\", \"93\u00d797=9021\"],\n    [\"34\u00d741=1394\", \"38\u00d729=1102\"],\n    [\"52\u00d720=1040\", \"88\u00d742=3696\"],\n    [\"58\u00d732=1856\", \"75\u00d755=4125\"],\n    [\"40\u00d739=1560\", \"90\u00d719=1710\"],\n    [\"66\u00d713=858\", \"64\u00d769=4416\"],\n    [\"29\u00d792=2668\", \"72\u00d718=1296\"],\n    [\"78\u00d798=7644\", \"88\u00d747=4136\"],\n    [\"40\u00d788=3520\", \"60\u00d731=1860\"],\n    [\"39\u00d797=3783\", \"52\u00d730=1560\"],\n    [\"74\u00d757=4218\", \"58\u00d730=1740\"],\n    [\"78\u00d728=2184\", \"51\u00d784=4284\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old = \"2023-09-03 Sunday\"; new = \"2023-09-04 Monday\"},\n    @{old = \"33\u00d747=1551\"; new = \"83\u00d711=913\"},\n    @{old = \"73\u00d768=4964\"; new = \"69\u00d714=966\"},\n    @{old = \"62\u00d742=2604\"; new = \"40\u00d751=2040\"},\n    @{old = \"79\u00d795=7505\"; new = \"25\u00d798=2450\"},\n    @{old = \"62\u00d747=2914\"; new = \"28\u00d780=2240\"},\n    @{old = \"72\u00d755=3960\"; new = \"41\u00d740=1640\"},\n    @{old = \"44\u00d728=1232\"; new = \"17\u00d713=221\"},\n    @{old = \"70\u00d723=1610\"; new = \"65\u00d758=3770\"},\n    @{old = \"43\u00d716=688\"; new = \"88\u00d726=2288\"},\n    @{old = \"42\u00d797=4074\"; new = \"92\u00d719=1748\"},\n    @{old = \"32\u00d796=3072\"; new = \"67\u00d725=1675\"},\n    @{old = \"19\u00d726=494\"; new = \"38\u00d731=1178\"},\n    @{old = \"27\u00d783=2241\"; new = \"98\u00d784=8232\"},\n    @{old = \"19\u00d756=1064\"; new = \"93\u00d797=9021\"},\n    @{old = \"34\u00d741=1394\"; new = \"38\u00d729=1102\"},\n    @{old = \"52\u00d720=1040\"; new = \"88\u00d742=3696\"},\n    @{old = \"58\u00d732=1856\"; new = \"75\u00d755=4125\"},\n    @{old = \"40\u00d739=1560\"; new = \"90\u00d719=1710\"},\n    @{old = \"66\u00d713=858\"; new = \"64\u00d769=4416\"},\n    @{old = \"29\u00d792=2668\"; new = \"72\u00d718=1296\"},\n    @{old = \"78\u00d798=7644\"; new = \"88\u00d747=4136\"},\n    @{old = \"40\u00d788=3520\"; new = \"60\u00d731=1860\"},\n    @{old = \"39\u00d797=3783\"; new = \"52\u00d730=1560\"},\n    @{old = \"74\u00d757=4218\"; new = \"58\u00d730=1740\"},\n    @{old = \"78\u00d728=2184\"; new = \"51\u00d784=4284\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.new\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n}\n"}
